# Applies the "Final data changes, finalize igor pro figures" edit:
#  - Scintillation Counter Results: append 12 new measurement rows (65-76)
#  - Count->Actual Activity: fill in C11:D13 (previously blank) which ripples
#    through Bottle Results and Averaged Results via existing formulas
#  - Update a few sheet selections to match where the author left the cursor

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Scintillation Counter Results (sheet2): new rows 65-76
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Scintillation Counter Results")
$ws2.Activate()

$ws2.Cells.Item(65,1).Value = 42993.404166666667
$ws2.Cells.Item(65,1).NumberFormat = "m/d/yy h:mm"
$ws2.Cells.Item(65,2).Value = "RaPYRASW_1A"
$ws2.Cells.Item(65,3).Value = 463.2
$ws2.Cells.Item(65,4).Value = 2.94
$ws2.Cells.Item(65,5).Value = 2.04
$ws2.Cells.Item(65,6).Value = 53.52
$ws2.Cells.Item(66,1).Value = 42993.404166666667
$ws2.Cells.Item(66,1).NumberFormat = "m/d/yy h:mm"
$ws2.Cells.Item(66,2).Value = "RaPYRASW_1B"
$ws2.Cells.Item(66,3).Value = 466.1
$ws2.Cells.Item(66,4).Value = 2.93
$ws2.Cells.Item(66,5).Value = 1.96
$ws2.Cells.Item(66,6).Value = 64.39
$ws2.Cells.Item(67,1).Value = 42993.404166666667
$ws2.Cells.Item(67,1).NumberFormat = "m/d/yy h:mm"
$ws2.Cells.Item(67,2).Value = "RaPYRASW_1C"
$ws2.Cells.Item(67,3).Value = 409
$ws2.Cells.Item(67,4).Value = 3.13
$ws2.Cells.Item(67,5).Value = 2.11
$ws2.Cells.Item(67,6).Value = 75.25
$ws2.Cells.Item(68,1).Value = 42993.554166666669
$ws2.Cells.Item(68,1).NumberFormat = "m/d/yy h:mm"
$ws2.Cells.Item(68,2).Value = "RaPYRASW_1A"
$ws2.Cells.Item(68,3).Value = 447
$ws2.Cells.Item(68,4).Value = 2.99
$ws2.Cells.Item(68,5).Value = 1.23
$ws2.Cells.Item(68,6).Value = 53.48
$ws2.Cells.Item(69,1).Value = 42993.554166666669
$ws2.Cells.Item(69,1).NumberFormat = "m/d/yy h:mm"
$ws2.Cells.Item(69,2).Value = "RaPYRASW_1B"
$ws2.Cells.Item(69,3).Value = 448.1
$ws2.Cells.Item(69,4).Value = 2.99
$ws2.Cells.Item(69,5).Value = 1.1910000000000001
$ws2.Cells.Item(69,6).Value = 64.209999999999994
$ws2.Cells.Item(70,1).Value = 42993.554166666669
$ws2.Cells.Item(70,1).NumberFormat = "m/d/yy h:mm"
$ws2.Cells.Item(70,2).Value = "RaPYRASW_1C"
$ws2.Cells.Item(70,3).Value = 422.1
$ws2.Cells.Item(70,4).Value = 3.08
$ws2.Cells.Item(70,5).Value = 1.19
$ws2.Cells.Item(70,6).Value = 74.959999999999994
$ws2.Cells.Item(71,1).Value = 42993.65625
$ws2.Cells.Item(71,1).NumberFormat = "m/d/yy h:mm"
$ws2.Cells.Item(71,2).Value = "RaPYRASW_1A"
$ws2.Cells.Item(71,3).Value = 446.7
$ws2.Cells.Item(71,4).Value = 2.99
$ws2.Cells.Item(71,5).Value = 0.7
$ws2.Cells.Item(71,6).Value = 53.43
$ws2.Cells.Item(72,1).Value = 42993.65625
$ws2.Cells.Item(72,1).NumberFormat = "m/d/yy h:mm"
$ws2.Cells.Item(72,2).Value = "RaPYRASW_1B"
$ws2.Cells.Item(72,3).Value = 438.4
$ws2.Cells.Item(72,4).Value = 3.02
$ws2.Cells.Item(72,5).Value = 0.71
$ws2.Cells.Item(72,6).Value = 64.14
$ws2.Cells.Item(73,1).Value = 42993.65625
$ws2.Cells.Item(73,1).NumberFormat = "m/d/yy h:mm"
$ws2.Cells.Item(73,2).Value = "RaPYRASW_1C"
$ws2.Cells.Item(73,3).Value = 428
$ws2.Cells.Item(73,4).Value = 3.06
$ws2.Cells.Item(73,5).Value = 0.68
$ws2.Cells.Item(73,6).Value = 74.95
$ws2.Cells.Item(74,1).Value = 42996.414583333331
$ws2.Cells.Item(74,1).NumberFormat = "m/d/yy h:mm"
$ws2.Cells.Item(74,2).Value = "RaPYRASW_1A"
$ws2.Cells.Item(74,3).Value = 442.6
$ws2.Cells.Item(74,4).Value = 3.01
$ws2.Cells.Item(74,5).Value = 0.07
$ws2.Cells.Item(74,6).Value = 10.63
$ws2.Cells.Item(75,1).Value = 42996.414583333331
$ws2.Cells.Item(75,1).NumberFormat = "m/d/yy h:mm"
$ws2.Cells.Item(75,2).Value = "RaPYRASW_1B"
$ws2.Cells.Item(75,3).Value = 426.9
$ws2.Cells.Item(75,4).Value = 3.06
$ws2.Cells.Item(75,5).Value = 0.08
$ws2.Cells.Item(75,6).Value = 21.28
$ws2.Cells.Item(76,1).Value = 42996.414583333331
$ws2.Cells.Item(76,1).NumberFormat = "m/d/yy h:mm"
$ws2.Cells.Item(76,2).Value = "RaPYRASW_1C"
$ws2.Cells.Item(76,3).Value = 431
$ws2.Cells.Item(76,4).Value = 3.05
$ws2.Cells.Item(76,5).Value = 0.08
$ws2.Cells.Item(76,6).Value = 31.91

$ws2.Range("G66").Select()

# ---------------------------------------------------------------------------
# Count->Actual Activity (sheet4): fill C11:D13
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Count->Actual Activity")
$ws4.Activate()

$ws4.Cells.Item(11,3).Value = 7.497916666666666
$ws4.Cells.Item(11,4).Value = 0.2236253645833333
$ws4.Cells.Item(12,3).Value = 7.414583333333334
$ws4.Cells.Item(12,4).Value = 0.2224375
$ws4.Cells.Item(13,3).Value = 7.042083333333333
$ws4.Cells.Item(13,4).Value = 0.2168961666666667

$ws4.Range("C11:D13").Select()

# ---------------------------------------------------------------------------
# Bottle Results (sheet5): selection only (values recompute automatically
# from the Count->Actual Activity formulas above)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Bottle Results")
$ws5.Activate()
$ws5.Range("W11").Select()

# ---------------------------------------------------------------------------
# Averaged Results (sheet6): selection only (values recompute automatically)
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Averaged Results")
$ws6.Activate()
$ws6.Range("B5:M5").Select()
